$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.720.47"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").Value = "1.652.70"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.0000"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.15"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3821"
$ws.Range("E7").Value = "  +2.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.38"
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3608"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.250"
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08242"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9998"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.63"
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.548"
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.404"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001234"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").Value = "1.652.06"
$ws.Range("E17").Value = "  +1.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.19"
$ws.Range("E18").Value = "  +3.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06970"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.787"
$ws.Range("E20").Value = "  +4.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.72"
$ws.Range("E21").Value = "  +0.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.63"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.568"
$ws.Range("E24").Value = "  +4.28%  "
$ws.Range("B25").Value = "WrappedBTC"
$ws.Range("C25").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D25").Value = "23.725.94"
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.074"
$ws.Range("E26").Value = "  -1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.34"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.24"
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.239"
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.25"
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").Value = "1.835.09"
$ws.Range("E31").Value = "  +1.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.879"
$ws.Range("E32").Value = "  +1.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.090"
$ws.Range("E33").Value = "  +5.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.95"
$ws.Range("E34").Value = "  +11.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.106"
$ws.Range("E35").Value = "  -5.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02831"
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2522"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08844"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.095"
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.07056"
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.83"
$ws.Range("E41").Value = "  +6.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7068"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.338"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.95"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6529"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.343"
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9996"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.984"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07990"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.15"
$ws.Range("E50").Value = "  +1.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.191"
$ws.Range("E51").Value = "  -0.32%  "
